$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in row 22 with a new logged session
$ws.Range("B22").Value = 45212
$ws.Range("C22").Value = 0.82291666666666663
$ws.Range("D22").Value = 0.90277777777777779
$ws.Range("G22").Value = "Got a lot of progress done on the file viewer/selector"
$ws.Range("H22").Value = "Next is making it actually open a file, also the terminal"

# Match row height auto-fit behavior seen in sibling rows
$ws.Rows.Item(22).AutoFit()

# Update the active selection like the author left it
$ws.Range("H23").Select()

$wb.Save()
